$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$strE2 = '[Enrico%Baldi%NULL%1,               Giuseppe Maria%Sechi%NULL%2,               Giuseppe Maria%Sechi%NULL%0,               Claudio%Mare%NULL%1,               Fabrizio%Canevari%NULL%1,               Antonella%Brancaglione%NULL%1,               Roberto%Primi%NULL%1,               Catherine%Klersy%NULL%1,               Alessandra%Palo%NULL%2,               Alessandra%Palo%NULL%0,               Enrico%Contri%NULL%1,               Vincenza%Ronchi%NULL%2,               Vincenza%Ronchi%NULL%0,               Giorgio%Beretta%NULL%1,               Francesca%Reali%NULL%1,               Pierpaolo%Parogni%NULL%1,               Fabio%Facchin%NULL%1,               Ugo%Rizzi%NULL%1,               Daniele%Bussi%NULL%1,               Simone%Ruggeri%NULL%1,               Luigi%Oltrona Visconti%NULL%1,               Simone%Savastano%s.savastano@smatteo.pv.it%1,               Sara%Compagnoni%NULL%3,               Sara%Compagnoni%NULL%0,               Sara%Compagnoni%NULL%0,               Rosa%Fracchia%NULL%1,               Antonio%Cuzzoli%NULL%1,               Andrea%Pagliosa%NULL%1,               Guido%Matiz%NULL%1,               Alessandra%Russo%NULL%1,               Andrea Lorenzo%Vecchi%NULL%1,               Cecilia%Fantoni%NULL%1,               Cristian%Fava%NULL%1,               Cinzia%Franzosi%NULL%1,               Claudio%Vimercati%NULL%1,               Dario%Franchi%NULL%1,               Enrico%Storti%NULL%1,               Erika%Taravelli%NULL%1,               Fulvio%Giovenzana%NULL%1,               Giovanni%Buetto%NULL%1,               Guido%Garzena%NULL%1,               Giorgio Antonio%Iotti%NULL%1,               Guido Francesco%Villa%NULL%1,               Marco%Botteri%NULL%1,               Salvatore Ivan%Caico%NULL%1,               Irene Raimondi%Cominesi%NULL%1,               Livio%Carnevale%NULL%1,               Matteo%Caresani%NULL%1,               Mario%Luppi%NULL%1,               Maurizio%Migliori%NULL%1,               Paola%Centineo%NULL%1,               Paola%Genoni%NULL%1,               Roberta%Bertona%NULL%1,               Roberto%De Ponti%NULL%1,               Riccardo%Osti%NULL%1,               Stefano%Buratti%NULL%1,               Gian Battista%Danzi%NULL%1,               Arianna%Marioni%NULL%1,               Antonella%De Pirro%NULL%1,               Simone%Molinari%NULL%1,               Vito%Sgromo%NULL%1,               Valeria%Musella%NULL%1,               Martina%Paglino%NULL%1,               Francesco%Mojoli%NULL%1,               Bruno%Lusona%NULL%1,               Michele%Pagani%NULL%1,               Moreno%Curti%NULL%1]'
$strE3 = '[Xiang%Wan%xwan@comp.hkbu.edu.hk%1,               Wenqian%Wang%wenqianwang2014@u.northwestern.edu%1,               Jiming%Liu%jiming@comp.hkbu.edu.hk%1,               Tiejun%Tong%tongt@hkbu.edu.hk%1]'
$strE4 = '[Pamela H.%Lai%NULL%1,               Elizabeth A.%Lancet%NULL%1,               Michael D.%Weiden%NULL%1,               Mayris P.%Webber%NULL%1,               Rachel%Zeig-Owens%NULL%1,               Charles B.%Hall%NULL%1,               David J.%Prezant%NULL%1]'
$strE5 = '[Andrea%Paoli%andrea.paoli@aopd.veneto.it%1,               Laura%Brischigliaro%NULL%2,               Laura%Brischigliaro%NULL%0,               Tommaso%Scquizzato%NULL%2,               Tommaso%Scquizzato%NULL%0,               Andrea%Favaretto%NULL%2,               Andrea%Favaretto%NULL%0,               Andrea%Spagna%NULL%1]'
$strE6 = '[Michael R.%Sayre%xref no email%1,        Leslie M.%Barnard%xref no email%1,        Catherine R.%Counts%xref no email%1,        Christopher J.%Drucker%xref no email%1,        Peter J.%Kudenchuk%xref no email%1,        Thomas D.%Rea%xref no email%1,        Mickey S.%Eisenberg%xref no email%1]'

$ws.Range("E2").Value = $strE2
$ws.Range("E3").Value = $strE3
$ws.Range("E4").Value = $strE4
$ws.Range("E5").Value = $strE5
$ws.Range("E6").Value = $strE6
